$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Separator column (F1:F5) holds the shared blank string " " -> " |"
# (points table index removed). Update every cell that uses this shared
# string so the string itself is edited in place rather than forked.
$ws.Range("F1").Value = " |"
$ws.Range("F2").Value = " |"
$ws.Range("F3").Value = " |"
$ws.Range("F4").Value = " |"
$ws.Range("F5").Value = " |"

# Group A table (columns B..E) and Group B table (columns H..K)
# Row 2
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 2

# Row 3
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 4
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 2

# Row 4
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 4
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6

# Row 5
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 2
